$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the distinct per-row IP addresses with a single shared "127.0.0.1" value
$ws.Range("F2").Value = "127.0.0.1"
$ws.Range("F3").Value = "127.0.0.1"
$ws.Range("F4").Value = "127.0.0.1"
$ws.Range("F5").Value = "127.0.0.1"
$ws.Range("F6").Value = "127.0.0.1"

# Update the active selection on the sheet
$ws.Range("F14").Select()
